$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

$ws.Range("C2").Value = 111.3
$ws.Range("C3").Value = 286.9
$ws.Range("C4").Value = 488.1
$ws.Range("C5").Value = 620.9
$ws.Range("C6").Value = 1271.6
$ws.Range("C7").Value = 1333
$ws.Range("C8").Value = 2573.7
$ws.Range("C9").Value = 2368.4

$ws.Range("F11").Select()

$wb.Save()
